# Generate two graphs for the website feature question type
#
# The "Table2" website-feature question (columns I, K, M) currently shares
# its answer-option labels with the "Table2" columns used for the other
# graph (columns J, L). Splitting the option codes lets the two graphs be
# generated independently:
#   "1,2,4;OO;2" -> "1,2,4;OC;2"
#   "1,4;OO;3"   -> "1,4;R;3"
# Only columns I, K and M (the first graph) are affected; J and L keep the
# original option codes used by the second graph.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targetColumns = @("I", "K", "M")
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($row = 2; $row -le $lastRow; $row++) {
    foreach ($col in $targetColumns) {
        $cell = $ws.Range($col + $row)
        $value = $cell.Value2

        if ($value -eq "1,2,4;OO;2") {
            $cell.Value2 = "1,2,4;OC;2"
        }
        elseif ($value -eq "1,4;OO;3") {
            $cell.Value2 = "1,4;R;3"
        }
    }
}
